$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.689801893040383
$ws.Range("D2").Value = -0.1977895957625399
$ws.Range("E2").Value = 0.00953015815020954
$ws.Range("F2").Value = 0.01584236461698259
$ws.Range("G2").Value = 0.01516645206850985
$ws.Range("H2").Value = -0.05184377338114116
$ws.Range("I2").Value = -0.1409981214786789
$ws.Range("J2").Value = 0.03196566811971526
$ws.Range("K2").Value = 0.0665734164232875
$ws.Range("L2").Value = 0.04435601442290548
$ws.Range("M2").Value = -0.01419361331937312
$ws.Range("N2").Value = -0.03413895618995193
$ws.Range("B3").Value = -0.689801893040383
$ws.Range("D3").Value = 0.1754789520186602
$ws.Range("E3").Value = 0.02771503092720217
$ws.Range("F3").Value = -0.05556955588644442
$ws.Range("G3").Value = -0.06463751991784057
$ws.Range("H3").Value = 0.1146910281989733
$ws.Range("I3").Value = 0.09652520864899328
$ws.Range("J3").Value = -0.05239032369778457
$ws.Range("K3").Value = -0.05126324650408413
$ws.Range("L3").Value = -0.003435706016208673
$ws.Range("M3").Value = 0.06539098660404966
$ws.Range("N3").Value = 0.05038073059792177
$ws.Range("B4").Value = -0.1977895957625399
$ws.Range("C4").Value = 0.1754789520186602
$ws.Range("E4").Value = 0.1837602300496707
$ws.Range("F4").Value = -0.1460803540235835
$ws.Range("G4").Value = -0.1345222434929064
$ws.Range("H4").Value = 0.3366930990988854
$ws.Range("I4").Value = -0.3482231603619768
$ws.Range("J4").Value = -0.1039164003795496
$ws.Range("K4").Value = -0.4081009320055506
$ws.Range("L4").Value = -0.3005519182748715
$ws.Range("M4").Value = 0.2776321694413547
$ws.Range("N4").Value = -0.3571763908955822
$ws.Range("B5").Value = 0.00953015815020954
$ws.Range("C5").Value = 0.02771503092720217
$ws.Range("D5").Value = 0.1837602300496707
$ws.Range("F5").Value = 0.3581846282573636
$ws.Range("G5").Value = 0.3614420468932035
$ws.Range("H5").Value = 0.4012543964666505
$ws.Range("I5").Value = 0.006822339822395155
$ws.Range("J5").Value = 0.7983442545732148
$ws.Range("K5").Value = 0.09569588899601536
$ws.Range("L5").Value = 0.01011896657585113
$ws.Range("M5").Value = 0.3944467797448968
$ws.Range("N5").Value = -0.01425156638275876
$ws.Range("B6").Value = 0.01584236461698259
$ws.Range("C6").Value = -0.05556955588644442
$ws.Range("D6").Value = -0.1460803540235835
$ws.Range("E6").Value = 0.3581846282573636
$ws.Range("G6").Value = 0.980928852548448
$ws.Range("H6").Value = 0.3685197946180032
$ws.Range("I6").Value = 0.1761391727802132
$ws.Range("J6").Value = 0.6185290620948303
$ws.Range("K6").Value = -0.04412919963610557
$ws.Range("L6").Value = -0.3118451534525333
$ws.Range("M6").Value = 0.3676699584686091
$ws.Range("N6").Value = 0.05412256302001409
$ws.Range("B7").Value = 0.01516645206850985
$ws.Range("C7").Value = -0.06463751991784057
$ws.Range("D7").Value = -0.1345222434929064
$ws.Range("E7").Value = 0.3614420468932035
$ws.Range("F7").Value = 0.980928852548448
$ws.Range("H7").Value = 0.3699904242157438
$ws.Range("I7").Value = 0.2355494801008387
$ws.Range("J7").Value = 0.6315276142144191
$ws.Range("K7").Value = 0.00934476116828869
$ws.Range("L7").Value = -0.2677596881566871
$ws.Range("M7").Value = 0.3769175723241422
$ws.Range("N7").Value = 0.07066182269695842
$ws.Range("B8").Value = -0.05184377338114116
$ws.Range("C8").Value = 0.1146910281989733
$ws.Range("D8").Value = 0.3366930990988854
$ws.Range("E8").Value = 0.4012543964666505
$ws.Range("F8").Value = 0.3685197946180032
$ws.Range("G8").Value = 0.3699904242157438
$ws.Range("I8").Value = -0.1889160726068343
$ws.Range("J8").Value = 0.4271370019900166
$ws.Range("K8").Value = -0.1407985109390405
$ws.Range("L8").Value = -0.2086571319298562
$ws.Range("M8").Value = 0.912656523852559
$ws.Range("N8").Value = -0.1531507197198079
$ws.Range("B9").Value = -0.1409981214786789
$ws.Range("C9").Value = 0.09652520864899328
$ws.Range("D9").Value = -0.3482231603619768
$ws.Range("E9").Value = 0.006822339822395155
$ws.Range("F9").Value = 0.1761391727802132
$ws.Range("G9").Value = 0.2355494801008387
$ws.Range("H9").Value = -0.1889160726068343
$ws.Range("J9").Value = 0.1892206722387704
$ws.Range("K9").Value = 0.6643566344192655
$ws.Range("L9").Value = 0.4960415350896452
$ws.Range("M9").Value = -0.2274648898758793
$ws.Range("N9").Value = 0.7496178209675342
$ws.Range("B10").Value = 0.03196566811971526
$ws.Range("C10").Value = -0.05239032369778457
$ws.Range("D10").Value = -0.1039164003795496
$ws.Range("E10").Value = 0.7983442545732148
$ws.Range("F10").Value = 0.6185290620948303
$ws.Range("G10").Value = 0.6315276142144191
$ws.Range("H10").Value = 0.4271370019900166
$ws.Range("I10").Value = 0.1892206722387704
$ws.Range("K10").Value = 0.2755031665414788
$ws.Range("L10").Value = 0.06563921175841496
$ws.Range("M10").Value = 0.4772401623332313
$ws.Range("N10").Value = 0.1723205790794998
$ws.Range("B11").Value = 0.0665734164232875
$ws.Range("C11").Value = -0.05126324650408413
$ws.Range("D11").Value = -0.4081009320055506
$ws.Range("E11").Value = 0.09569588899601536
$ws.Range("F11").Value = -0.04412919963610557
$ws.Range("G11").Value = 0.00934476116828869
$ws.Range("H11").Value = -0.1407985109390405
$ws.Range("I11").Value = 0.6643566344192655
$ws.Range("J11").Value = 0.2755031665414788
$ws.Range("L11").Value = 0.9291646513948305
$ws.Range("M11").Value = -0.1302229061075761
$ws.Range("N11").Value = 0.7526418513922563
$ws.Range("B12").Value = 0.04435601442290548
$ws.Range("C12").Value = -0.003435706016208673
$ws.Range("D12").Value = -0.3005519182748715
$ws.Range("E12").Value = 0.01011896657585113
$ws.Range("F12").Value = -0.3118451534525333
$ws.Range("G12").Value = -0.2677596881566871
$ws.Range("H12").Value = -0.2086571319298562
$ws.Range("I12").Value = 0.4960415350896452
$ws.Range("J12").Value = 0.06563921175841496
$ws.Range("K12").Value = 0.9291646513948305
$ws.Range("M12").Value = -0.2054599720196534
$ws.Range("N12").Value = 0.6915138823168836
$ws.Range("B13").Value = -0.01419361331937312
$ws.Range("C13").Value = 0.06539098660404966
$ws.Range("D13").Value = 0.2776321694413547
$ws.Range("E13").Value = 0.3944467797448968
$ws.Range("F13").Value = 0.3676699584686091
$ws.Range("G13").Value = 0.3769175723241422
$ws.Range("H13").Value = 0.912656523852559
$ws.Range("I13").Value = -0.2274648898758793
$ws.Range("J13").Value = 0.4772401623332313
$ws.Range("K13").Value = -0.1302229061075761
$ws.Range("L13").Value = -0.2054599720196534
$ws.Range("N13").Value = -0.1903295950547503
$ws.Range("B14").Value = -0.03413895618995193
$ws.Range("C14").Value = 0.05038073059792177
$ws.Range("D14").Value = -0.3571763908955822
$ws.Range("E14").Value = -0.01425156638275876
$ws.Range("F14").Value = 0.05412256302001409
$ws.Range("G14").Value = 0.07066182269695842
$ws.Range("H14").Value = -0.1531507197198079
$ws.Range("I14").Value = 0.7496178209675342
$ws.Range("J14").Value = 0.1723205790794998
$ws.Range("K14").Value = 0.7526418513922563
$ws.Range("L14").Value = 0.6915138823168836
$ws.Range("M14").Value = -0.1903295950547503
